# Update the quarterly database: drop the oldest quarter column and
# append the newest quarter ("فصل دوم منتهی به 1401/12"), shifting every
# quarter-indexed value one column to the left (read_price algorithm).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues($range, $values) {
    $arr = New-Object 'object[,]' 1,$values.Length
    for ($i = 0; $i -lt $values.Length; $i++) {
        $arr[0,$i] = $values[$i]
    }
    $ws.Range($range).Value = $arr
}

# --- Header rows: the 10 rolling quarter labels (columns E:N) ---
$quarterHeaders = @(
    "فصل اول منتهی به 1399/09",
    "فصل دوم منتهی به 1399/12",
    "فصل سوم منتهی به 1400/03",
    "فصل چهارم منتهی به 1400/06",
    "فصل اول منتهی به 1400/09",
    "فصل دوم منتهی به 1400/12",
    "فصل سوم منتهی به 1401/03",
    "فصل چهارم منتهی به 1401/06",
    "فصل اول منتهی به 1401/09",
    "فصل دوم منتهی به 1401/12"
)

Set-RowValues "E8:N8" $quarterHeaders
Set-RowValues "E24:N24" $quarterHeaders

# --- Data rows: shift existing quarterly figures left by one column and
#     bring in the newly reported quarter's figure in column N ---
Set-RowValues "E11:N11" @(0, 0, 0, 0, 0, 0, 0, 81102, 16233, 40304)
Set-RowValues "E13:N13" @(260, 2488, -2748, 36484, 7150, 2304, 5915, 50419, 19840, 79423)
Set-RowValues "E16:N16" @(705, 121, 924, 1284, 942, 1021, 1012, 967, 1112, 1091)
Set-RowValues "E17:N17" @(2004, 8527, 4242, 16190, 9385, 4433, 543, 28118, 14643, 18667)
Set-RowValues "E18:N18" @(0, 0, 0, 57230, 0, 0, 0, 0, 0, 0)
Set-RowValues "E19:N19" @(16211, 12010, 23684, 54943, 19497, 33959, 80116, 71505, 66715, 147343)
Set-RowValues "E20:N20" @(19180, 23146, 26102, 166131, 36974, 41717, 87586, 232111, 118543, 286828)
Set-RowValues "E26:N26" @(168, 162, 175, 210, 206, 165, 148, 208, 213, 210)
Set-RowValues "E27:N27" @(755, 756, 765, 691, 715, 740, 764, 697, 697, 706)
